$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "A24"
$ws.Range("D6").Value = "B24"
$ws.Range("D7").Value = "C24"
$ws.Range("D8").Value = "G24"
$ws.Range("D9").Value = "H24"
$ws.Range("D10").Value = "I24"
$ws.Range("D11").Value = "J24"

$ws.Range("D5:D11").Select()
